$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update columns G2:AV2 (42 cells)
$row2 = New-Object "object[,]" 1,42
$row2[0,0] = 172.825215362831
$row2[0,1] = 175.416851512284
$row2[0,2] = 177.199196639361
$row2[0,3] = 178.999651509742
$row2[0,4] = 180.818400129766
$row2[0,5] = 182.655628375392
$row2[0,6] = 243.397645363672
$row2[0,7] = 245.827873526138
$row2[0,8] = 248.28236654504
$row2[0,9] = 250.761366695267
$row2[0,10] = 252.231435789332
$row2[0,11] = 253.710123049624
$row2[0,12] = 0
$row2[0,13] = 0
$row2[0,14] = 0
$row2[0,15] = 0
$row2[0,16] = 0
$row2[0,17] = 0
$row2[0,18] = 0
$row2[0,19] = 0
$row2[0,20] = 0
$row2[0,21] = 0
$row2[0,22] = 0
$row2[0,23] = 0
$row2[0,24] = 0
$row2[0,25] = 0
$row2[0,26] = 0
$row2[0,27] = 0
$row2[0,28] = 0
$row2[0,29] = 0
$row2[0,30] = 0
$row2[0,31] = 0
$row2[0,32] = 0
$row2[0,33] = 0
$row2[0,34] = 0
$row2[0,35] = 0
$row2[0,36] = 0
$row2[0,37] = 0
$row2[0,38] = 0
$row2[0,39] = 0
$row2[0,40] = 0
$row2[0,41] = 0
$ws.Range("G2:AV2").Value = $row2

# Row 3: update columns G3:AV3 (42 cells)
$row3 = New-Object "object[,]" 1,42
$row3[0,0] = 159.780785214987
$row3[0,1] = 170.057946057514
$row3[0,2] = 176.793058737561
$row3[0,3] = 183.676060335245
$row3[0,4] = 190.709714324339
$row3[0,5] = 196.79829535043
$row3[0,6] = 202.020669289134
$row3[0,7] = 206.201363920483
$row3[0,8] = 210.450353309583
$row3[0,9] = 214.613968850182
$row3[0,10] = 222.279856734501
$row3[0,11] = 230.338652341194
$row3[0,12] = 238.752088321691
$row3[0,13] = 246.605298285589
$row3[0,14] = 247.052953509889
$row3[0,15] = 247.29665334348
$row3[0,16] = 247.333986338494
$row3[0,17] = 247.200618945643
$row3[0,18] = 247.604756422399
$row3[0,19] = 247.853253102781
$row3[0,20] = 247.982091943821
$row3[0,21] = 247.951000827654
$row3[0,22] = 246.788242201678
$row3[0,23] = 245.432408562552
$row3[0,24] = 243.918748838027
$row3[0,25] = 242.244356508887
$row3[0,26] = 241.038171302294
$row3[0,27] = 239.646592519882
$row3[0,28] = 238.105116499942
$row3[0,29] = 236.372976725091
$row3[0,30] = 234.799177029243
$row3[0,31] = 233.038045405906
$row3[0,32] = 231.086997551258
$row3[0,33] = 228.981526472646
$row3[0,34] = 226.773609543596
$row3[0,35] = 224.406376101317
$row3[0,36] = 221.915218642343
$row3[0,37] = 219.297379214267
$row3[0,38] = 321.720903844883
$row3[0,39] = 416.839625234504
$row3[0,40] = 505.29673765003
$row3[0,41] = 587.57149317354
$ws.Range("G3:AV3").Value = $row3

# Row 6: update columns G6:AV6 (42 cells)
$row6 = New-Object "object[,]" 1,42
$row6[0,0] = 1386.21726265027
$row6[0,1] = 1368.22675839617
$row6[0,2] = 1401.99671443597
$row6[0,3] = 1435.8215732094
$row6[0,4] = 1469.70189256263
$row6[0,5] = 1503.6382360099
$row6[0,6] = 1469.10158396751
$row6[0,7] = 1492.95013480382
$row6[0,8] = 1516.84856784654
$row6[0,9] = 1540.79738114956
$row6[0,10] = 1551.2342805167
$row6[0,11] = 1561.68889658237
$row6[0,12] = 0
$row6[0,13] = 0
$row6[0,14] = 0
$row6[0,15] = 0
$row6[0,16] = 0
$row6[0,17] = 0
$row6[0,18] = 0
$row6[0,19] = 0
$row6[0,20] = 0
$row6[0,21] = 0
$row6[0,22] = 0
$row6[0,23] = 0
$row6[0,24] = 0
$row6[0,25] = 0
$row6[0,26] = 0
$row6[0,27] = 0
$row6[0,28] = 0
$row6[0,29] = 0
$row6[0,30] = 0
$row6[0,31] = 0
$row6[0,32] = 0
$row6[0,33] = 0
$row6[0,34] = 0
$row6[0,35] = 0
$row6[0,36] = 0
$row6[0,37] = 0
$row6[0,38] = 0
$row6[0,39] = 0
$row6[0,40] = 0
$row6[0,41] = 0
$ws.Range("G6:AV6").Value = $row6
